$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet ---
# URL changed from national-directory-query to fhir-directory-query
$wsMeta.Range("B2").Value = "http://hl7.org/fhir/us/fhir-directory-query/StructureDefinition/accessibility"

# Date updated to reflect the new publication run
$wsMeta.Range("B8").Value = "2021-12-17T13:53:37-05:00"

# --- Elements sheet ---
# Same StructureDefinition URL repeated in the elements table (Base Definition-ish column)
$wsElem.Range("Q5").Value = "http://hl7.org/fhir/us/fhir-directory-query/StructureDefinition/accessibility"

# Binding Value Set URL for the Accessibility value set
$wsElem.Range("Y6").Value = "http://hl7.org/fhir/us/fhir-directory-query/ValueSet/AccessibilityVS"

# The shorter URL text narrows the best-fit width of column Y
$wsElem.Columns.Item(25).ColumnWidth = 60.67
